$d = $word.ActiveDocument

# wdBrightGreen (WdColorIndex = 4) serializes to OOXML <w:highlight w:val="green"/>
$wdBrightGreen = 4

# The three bullet paragraphs that got green-highlighted in the commit:
#   "Brand (ishonch, texnomart, mediapark) uchun CRUD. ..."
#   "Branch (brandlarni filiallari) uchun CRUD. ..."
#   "Regionlar kesimida, Tumanlardagi har bir brandning branchlar sonini olish imkoniyati."
# Identify them by their distinctive leading text so the script is resilient
# to the exact paragraph index.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Brand (") -or $t.StartsWith("Branch (") -or $t.StartsWith("Regionlar ")) {
        # Setting Font.HighlightColorIndex on the whole paragraph range (which
        # includes the trailing paragraph mark) highlights every run plus the
        # paragraph mark itself, matching the target markup exactly.
        $p.Range.Font.HighlightColorIndex = $wdBrightGreen
    }
}
